{"js": "// Fix typo: \"Moskou laten zien\" -> \"Moscow laten zien\"\nconst results = context.document.body.search(\"Moskou laten zien\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Moscow laten zien\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Fix typo: \"Moskou laten zien\" -> \"Moscow laten zien\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"Moskou laten zien\", $false, $false, $false, $false, $false, $true, 1, $false, \"Moscow laten zien\", 2) | Out-Null\n"}
